$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 7819
$ws.Range("F5").Value = 7819
$ws.Range("F14").Value = 5745
$ws.Range("F16").Value = 2719
$ws.Range("F23").Value = 41
$ws.Range("F24").Value = 3841
$ws.Range("F26").Value = 56
$ws.Range("F27").Value = 49
$ws.Range("F30").Value = 5340
$ws.Range("F34").Value = 385
$ws.Range("F37").Value = 1810
$ws.Range("F38").Value = 1005
$ws.Range("F41").Value = 3821
$ws.Range("F45").Value = 3463

# Sheet: 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 15
$ws.Range("F6").Value = 15

# Sheet: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1359

# Sheet: 全部类型 (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1359
$ws.Range("F5").Value = 7819
$ws.Range("F6").Value = 7819
$ws.Range("F13").Value = 5745
$ws.Range("F15").Value = 2719
$ws.Range("F24").Value = 15
$ws.Range("F25").Value = 3841
$ws.Range("F27").Value = 56
$ws.Range("F28").Value = 49
$ws.Range("F31").Value = 5340
$ws.Range("F34").Value = 385
$ws.Range("F38").Value = 1811
$ws.Range("F39").Value = 1005
$ws.Range("F43").Value = 3821
$ws.Range("F47").Value = 3463

$wb.Save()
